$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.118266759376858
$ws.Range("C2").Value = -0.0500365083534178
$ws.Range("D2").Value = 0.0898747801425915
$ws.Range("E2").Value = 0.773046940540996
$ws.Range("F2").Value = 0.169379849461424
$ws.Range("G2").Value = -0.398910531419891
$ws.Range("H2").Value = -0.323633203048968
$ws.Range("I2").Value = 0.0447560739024471
$ws.Range("J2").Value = -0.129806726007518
$ws.Range("K2").Value = 0.570091796525703
$ws.Range("L2").Value = 0.703817607112367
$ws.Range("M2").Value = -0.108130951189116
$ws.Range("N2").Value = 0.538430729599699
$ws.Range("O2").Value = -0.545033936855826
